$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 272 (pushes old rows 272..351 down to 274..353),
# mirroring a new week of price data being added to the top of this rotating 80-week window
# and the oldest week's rows being carried along at the tail.
$ws.Range("A272:A273").EntireRow.Insert()

# Row 272: Acelga, Primera, new week (2022-03-07 / serial 44627)
$ws.Cells.Item(272, 1).Value = 8
$ws.Cells.Item(272, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(272, 3).Value = "Coquimbo"
$ws.Cells.Item(272, 4).Value = 44627
$ws.Cells.Item(272, 5).Value = 4
$ws.Cells.Item(272, 6).Value = 100112009
$ws.Cells.Item(272, 7).Value = "Acelga"
$ws.Cells.Item(272, 8).Value = "Sin especificar"
$ws.Cells.Item(272, 9).Value = "Primera"
$ws.Cells.Item(272, 10).Value = 2500
$ws.Cells.Item(272, 11).Value = 500
$ws.Cells.Item(272, 12).Value = 600
$ws.Cells.Item(272, 13).Value = 550
$ws.Cells.Item(272, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(272, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(272, 16).Value = 275
$ws.Cells.Item(272, 17).Value = 2
$ws.Cells.Item(272, 18).Value = "Hortaliza"

# Row 273: Acelga, Segunda, same new week
$ws.Cells.Item(273, 1).Value = 8
$ws.Cells.Item(273, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(273, 3).Value = "Coquimbo"
$ws.Cells.Item(273, 4).Value = 44627
$ws.Cells.Item(273, 5).Value = 4
$ws.Cells.Item(273, 6).Value = 100112009
$ws.Cells.Item(273, 7).Value = "Acelga"
$ws.Cells.Item(273, 8).Value = "Sin especificar"
$ws.Cells.Item(273, 9).Value = "Segunda"
$ws.Cells.Item(273, 10).Value = 1320
$ws.Cells.Item(273, 11).Value = 400
$ws.Cells.Item(273, 12).Value = 450
$ws.Cells.Item(273, 13).Value = 425
$ws.Cells.Item(273, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(273, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(273, 16).Value = 212
$ws.Cells.Item(273, 17).Value = 2
$ws.Cells.Item(273, 18).Value = "Hortaliza"

# Match the date-column number format on the two new rows to the rest of column D
$ws.Cells.Item(272, 4).NumberFormat = $ws.Cells.Item(274, 4).NumberFormat
$ws.Cells.Item(273, 4).NumberFormat = $ws.Cells.Item(274, 4).NumberFormat
